# Storage bulk-upload template: add a "Patient ID" column after "Sample Code"
# and drop the leftover example data rows (2-6), leaving only the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the sample data rows (rows 2-6, column A VL-codes) - only the
#    header row should remain.
$ws.Range("A2:A6").EntireRow.Delete()

# 2) Insert a new column before column B ("Feezer Code" and everything to its
#    right shifts one column to the right). The new column inherits the
#    header formatting (red bold) from column A automatically.
$ws.Range("B1").EntireColumn.Insert()

# 3) Give the new column a sensible width (approximating 21.95 characters;
#    the COM layer rounds to whole pixels, same as Excel itself does).
$ws.Columns.Item(2).ColumnWidth = 21.1666666666667

# 4) Label the new header cell and restyle it to match the other "blank"
#    header cells (black bold) instead of the red bold used for labelled
#    headers - copy the look from H1 (an already black-bold header cell)
#    and then set the text.
$ws.Range("H1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B1").Value2 = "Patient ID"

# Note: columns 3-12 (formerly 2-11) keep their original widths automatically
# because EntireColumn.Insert() shifts the pre-existing per-column formatting
# along with them - no need to (and better not to) touch them again here.

# 5) The "Volume (ml)" label moved from column F to column G after the
#    insert; give it the labelled-header look (red bold) to match the rest
#    of the labelled headers.
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G1").Value2 = "Volume (ml)"
